# Daily attendance processing - 2026-02-01 18:48:03
# Swap the order of the "Recorded By" entries in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$xlWhole = 1
$xlByRows = 1

$colG = $ws.Columns.Item(7)
$colG.Replace($oldValue, $newValue, $xlWhole, $xlByRows, $false, $false, $false) | Out-Null
